$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.247.96'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.862.74'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.81%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9984'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7153'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '240.71'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9996'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3088'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07721'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.00'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08318'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.879.52'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7187'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.219'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '90.94'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.254.19'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.990'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.93%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '243.90'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.149.47'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007814'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.18'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9991'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.952'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9994'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1616'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.82'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.917'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.61'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.358'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.53%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.445'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.01%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.497'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.258'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05188'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.8187'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +13.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.935'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.176'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.678'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01860'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('E39').Value = '  -1.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.162.97'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.213'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9009'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '72.94'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.32%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9990'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('E45').Value = '  -1.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.044.40'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.41%  '
$ws.Range('E47').Value = '  -2.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.786'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.384'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.29%  '
$ws.Range('E50').Value = '  -0.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.085'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.03%  '
